# Applies the "Updated symbol list" commit: refreshes cryptocurrency prices
# and shifts the Coin/Link/Volume columns for rows 18-24 down by one entry
# (a new coin, "One", was inserted at the top of that block while the
# remaining coins moved down one slot), keeping each row's own freshly
# scraped Price value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price (column D) refreshes -------------------------------------------------
$ws.Range("D2").Value  = "229.31"
$ws.Range("D3").Value  = "22.53"
$ws.Range("D4").Value  = "5.269"
$ws.Range("D5").Value  = "0.05570"
$ws.Range("D6").Value  = "3.380"
$ws.Range("D7").Value  = "6.472"
$ws.Range("D8").Value  = "1.048"
$ws.Range("D9").Value  = "0.7826"
$ws.Range("D11").Value = "0.07353"
$ws.Range("D12").Value = "0.03162"
$ws.Range("D14").Value = "0.09275"
$ws.Range("D15").Value = "0.001664"
$ws.Range("D16").Value = "3.266"
$ws.Range("D17").Value = "0.04786"

# --- Rows 18-24: Coin / Link / Volume shift down one slot, Price refreshed -------------
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005896"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "0.006242"
$ws.Range("E19").Value = "18TigerCashTCH"

$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "0.005232"
$ws.Range("E20").Value = "19HotbitTokenHTB"

$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "0.001057"
$ws.Range("E21").Value = "20BitKanKAN"

$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "0.0001502"
$ws.Range("E22").Value = "21NitroExNTX"

$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "3.939"
$ws.Range("E23").Value = "22LEOLEO"

$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "2.146"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# --- Remaining scattered updates --------------------------------------------------------
$ws.Range("D26").Value = "0.1308"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"

$ws.Range("D40").Value = "0.04000"
$ws.Range("D41").Value = "0.007042"
$ws.Range("D42").Value = "0.003388"
$ws.Range("D43").Value = "0.1038"
$ws.Range("D44").Value = "0.009999"
$ws.Range("D48").Value = "0.04231"
